$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Howard Livingston-" + "Github" + " Links"  ->  single run
#    "Howard Livingston-Github Links" (also drops the spellStart/
#    spellEnd proofErr markers that wrapped "Github").
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("Howard Livingston-Github Links", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Howard Livingston-Github Links", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "10/12" -> "11" / bookmark "_GoBack" / "/12"
#    First bump the text, then relocate the _GoBack bookmark so it
#    sits between "11" and "/12" (splitting the run in two), and
#    remove it from wherever it used to live (last paragraph).
# ------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4.Range.Find.Execute("10/12", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "11/12", 2) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p4b = $d.Paragraphs(4)
$r4 = $p4b.Range
$splitPoint = $d.Range($r4.Start + 2, $r4.Start + 2)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# ------------------------------------------------------------------
# 3) "Gitlhub" + " Link Live app" -> single run "Gitlhub Link Live app"
#    (drops the spellStart/spellEnd proofErr markers around "Gitlhub").
#    The proofErr sits as the very first child of the paragraph, so a
#    find/replace that starts exactly there leaves a stray marker
#    behind; widen the search to include the previous (empty)
#    paragraph's mark so the match does not start on the marker, and
#    re-insert that paragraph mark via ^p in the replacement text.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Gitlhub Link Live app*") {
        $gitlhubIdx = $i
        break
    }
}
$rangeWide = $d.Range($d.Paragraphs($gitlhubIdx - 1).Range.Start, $d.Paragraphs($gitlhubIdx).Range.End)
$rangeWide.MoveEnd(1, -1) | Out-Null
$rangeWide.Find.Execute($rangeWide.Text, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "^pGitlhub Link Live app", 2) | Out-Null

# ------------------------------------------------------------------
# 4) "Github" + " Master branch" -> single run "Github Master branch"
#    Same proofErr-at-paragraph-start situation as above.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Github Master branch*") {
        $masterIdx = $i
        break
    }
}
$rangeWide2 = $d.Range($d.Paragraphs($masterIdx - 1).Range.Start, $d.Paragraphs($masterIdx).Range.End)
$rangeWide2.MoveEnd(1, -1) | Out-Null
$rangeWide2.Find.Execute($rangeWide2.Text, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "^pGithub Master branch", 2) | Out-Null
